$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.Value = "'56.453.18"
$c.ClearFormats()

$c = $ws.Range("D3")
$c.Value = "'3.253.36"
$c.ClearFormats()

$c = $ws.Range("E3")
$c.Value = "'  +6.45%  "
$c.ClearFormats()

$c = $ws.Range("E4")
$c.Value = "'  +0.05%  "
$c.ClearFormats()

$c = $ws.Range("D5")
$c.Value = "'398.74"
$c.ClearFormats()

$c = $ws.Range("E5")
$c.Value = "'  +2.23%  "
$c.ClearFormats()

$c = $ws.Range("D6")
$c.Value = "'111.09"
$c.ClearFormats()

$c = $ws.Range("E6")
$c.Value = "'  +10.13%  "
$c.ClearFormats()

$c = $ws.Range("E7")
$c.Value = "'  +5.09%  "
$c.ClearFormats()

$c = $ws.Range("E8")
$c.Value = "'  -0.03%  "
$c.ClearFormats()

$c = $ws.Range("E9")
$c.Value = "'  +6.93%  "
$c.ClearFormats()

$c = $ws.Range("D10")
$c.Value = "'39.50"
$c.ClearFormats()

$c = $ws.Range("E10")
$c.Value = "'  +7.69%  "
$c.ClearFormats()

$c = $ws.Range("E11")
$c.Value = "'  +11.67%  "
$c.ClearFormats()

$c = $ws.Range("E12")
$c.Value = "'  +2.42%  "
$c.ClearFormats()

$c = $ws.Range("D13")
$c.Value = "'3.768.41"
$c.ClearFormats()

$c = $ws.Range("E13")
$c.Value = "'  +6.59%  "
$c.ClearFormats()

$c = $ws.Range("D14")
$c.Value = "'19.24"
$c.ClearFormats()

$c = $ws.Range("E14")
$c.Value = "'  +5.55%  "
$c.ClearFormats()

$c = $ws.Range("D15")
$c.Value = "'8.09"
$c.ClearFormats()

$c = $ws.Range("E15")
$c.Value = "'  +5.82%  "
$c.ClearFormats()

$c = $ws.Range("D16")
$c.Value = "'3.248.32"
$c.ClearFormats()

$c = $ws.Range("E16")
$c.Value = "'  +6.16%  "
$c.ClearFormats()

$c = $ws.Range("E17")
$c.Value = "'  +5.45%  "
$c.ClearFormats()

$c = $ws.Range("D18")
$c.Value = "'10.98"
$c.ClearFormats()

$c = $ws.Range("E18")
$c.Value = "'  +3.93%  "
$c.ClearFormats()

$c = $ws.Range("D19")
$c.Value = "'56.391.91"
$c.ClearFormats()

$c = $ws.Range("E19")
$c.Value = "'  +10.35%  "
$c.ClearFormats()

$c = $ws.Range("E20")
$c.Value = "'  +5.14%  "
$c.ClearFormats()

$c = $ws.Range("E21")
$c.Value = "'  +8.54%  "
$c.ClearFormats()

$c = $ws.Range("E22")
$c.Value = "'  +6.57%  "
$c.ClearFormats()

$c = $ws.Range("D23")
$c.Value = "'299.14"
$c.ClearFormats()

$c = $ws.Range("E23")
$c.Value = "'  +13.56%  "
$c.ClearFormats()

$c = $ws.Range("D24")
$c.Value = "'74.95"
$c.ClearFormats()

$c = $ws.Range("E24")
$c.Value = "'  +7.81%  "
$c.ClearFormats()

$c = $ws.Range("E25")
$c.Value = "'  +1.88%  "
$c.ClearFormats()

$c = $ws.Range("D26")
$c.Value = "'8.10"
$c.ClearFormats()

$c = $ws.Range("E26")
$c.Value = "'  +3.00%  "
$c.ClearFormats()

$c = $ws.Range("D27")
$c.Value = "'28.14"
$c.ClearFormats()

$c = $ws.Range("E27")
$c.Value = "'  +5.28%  "
$c.ClearFormats()

$c = $ws.Range("E28")
$c.Value = "'  +5.09%  "
$c.ClearFormats()

$c = $ws.Range("E29")
$c.Value = "'  +3.34%  "
$c.ClearFormats()

$c = $ws.Range("E30")
$c.Value = "'  +4.52%  "
$c.ClearFormats()

$c = $ws.Range("D31")
$c.Value = "'1.00"
$c.ClearFormats()

$c = $ws.Range("E31")
$c.Value = "'  +0.03%  "
$c.ClearFormats()

$c = $ws.Range("E32")
$c.Value = "'  +6.34%  "
$c.ClearFormats()

$c = $ws.Range("D33")
$c.Value = "'11.09"
$c.ClearFormats()

$c = $ws.Range("E33")
$c.Value = "'  +5.83%  "
$c.ClearFormats()

$c = $ws.Range("D34")
$c.Value = "'38.42"
$c.ClearFormats()

$c = $ws.Range("E34")
$c.Value = "'  +7.34%  "
$c.ClearFormats()

$c = $ws.Range("D35")
$c.Value = "'0.0488"
$c.ClearFormats()

$c = $ws.Range("E35")
$c.Value = "'  -0.39%  "
$c.ClearFormats()

$c = $ws.Range("E36")
$c.Value = "'  +5.66%  "
$c.ClearFormats()

$c = $ws.Range("D37")
$c.Value = "'51.56"
$c.ClearFormats()

$c = $ws.Range("E37")
$c.Value = "'  +3.16%  "
$c.ClearFormats()

$c = $ws.Range("E38")
$c.Value = "'  +28.40%  "
$c.ClearFormats()

$c = $ws.Range("E39")
$c.Value = "'  +5.35%  "
$c.ClearFormats()

$c = $ws.Range("D40")
$c.Value = "'0.999"
$c.ClearFormats()

$c = $ws.Range("E40")
$c.Value = "'  -0.10%  "
$c.ClearFormats()

$c = $ws.Range("D41")
$c.Value = "'17.57"
$c.ClearFormats()

$c = $ws.Range("E41")
$c.Value = "'  +6.36%  "
$c.ClearFormats()

$c = $ws.Range("E42")
$c.Value = "'  +6.53%  "
$c.ClearFormats()

$c = $ws.Range("D43")
$c.Value = "'133.86"
$c.ClearFormats()

$c = $ws.Range("E43")
$c.Value = "'  +3.21%  "
$c.ClearFormats()

$c = $ws.Range("E44")
$c.Value = "'  +4.64%  "
$c.ClearFormats()

$c = $ws.Range("E45")
$c.Value = "'  +6.21%  "
$c.ClearFormats()

$c = $ws.Range("E46")
$c.Value = "'  -2.37%  "
$c.ClearFormats()

$c = $ws.Range("D47")
$c.Value = "'22.20"
$c.ClearFormats()

$c = $ws.Range("E47")
$c.Value = "'  +2.51%  "
$c.ClearFormats()

$c = $ws.Range("D48")
$c.Value = "'2.149.81"
$c.ClearFormats()

$c = $ws.Range("E48")
$c.Value = "'  +4.16%  "
$c.ClearFormats()

$c = $ws.Range("E49")
$c.Value = "'  +1.58%  "
$c.ClearFormats()

$c = $ws.Range("D50")
$c.Value = "'2.42"
$c.ClearFormats()

$c = $ws.Range("E50")
$c.Value = "'  -2.23%  "
$c.ClearFormats()

$c = $ws.Range("D51")
$c.Value = "'1.98"
$c.ClearFormats()

$c = $ws.Range("E51")
$c.Value = "'  +39.52%  "
$c.ClearFormats()
